$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells hold text values (coin prices/percentages stored as strings,
# some of which look numeric, e.g. "0.0000206" or "2.50"). Force the cell format
# to Text before assignment so Excel does not silently convert them to numbers
# (which would also strip significant trailing zeros).
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '66.199.30'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -0.09%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.564.27'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +0.73%  '
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.11%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '605.49'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '144.19'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -0.25%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '3.563.54'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.73%  '
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +0.17%  '
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +2.18%  '
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -0.45%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '7.77'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  -3.45%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.413'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -0.33%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.170.27'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +0.77%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.0000206'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -1.34%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '30.33'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -0.49%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.570.53'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +0.97%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '66.278.39'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -0.14%  '
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -0.58%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.41'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +4.56%  '
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -0.34%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '14.78'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -1.63%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '430.39'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +0.98%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.613'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +1.74%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '79.53'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +1.26%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '3.705.63'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +0.60%  '
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -2.91%  '
$c = $ws.Range('B28')
$c.NumberFormat = '@'
$c.Value = 'PancakeSwap'
$c = $ws.Range('C28')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.50'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +0.74%  '
$c = $ws.Range('B29')
$c.NumberFormat = '@'
$c.Value = 'InternetComputer(DFINITY)'
$c = $ws.Range('C29')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.15'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -1.98%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '7.91'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -1.52%  '
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.559.98'
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +0.89%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '25.43'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +0.16%  '
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -2.84%  '
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -3.66%  '
$c = $ws.Range('B36')
$c.NumberFormat = '@'
$c.Value = 'USDe'
$c = $ws.Range('C36')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c = $ws.Range('B37')
$c.NumberFormat = '@'
$c.Value = 'Aptos'
$c = $ws.Range('C37')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '7.82'
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -0.67%  '
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -1.97%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '5.60'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -0.63%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '175.19'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +3.15%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0849'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -1.56%  '
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +0.13%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.887'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -0.65%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.92'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +1.62%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '45.99'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +1.02%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +2.01%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.19'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -2.02%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '24.90'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -4.56%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '7.13'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -1.01%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '23.36'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +3.82%  '
